# B6-PowerPoint.pptx edit
#  1. Re-style the three data tables (slides 14, 15, 16) from the deck's
#     custom table style onto the built-in "No Style, Table Grid" style.
#  2. Re-colour the presentation's theme palette from the "Integral /
#     Red Violet" scheme over to the stock "Office" scheme.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Tables: switch styleId on every table in the deck.
# ---------------------------------------------------------------------
$oldStyleId = "{63662A86-823A-481B-AD0B-CAA589A7E9DD}"
$newStyleId = "{5EDE5D6C-459C-4361-A20C-DFBF5E0A35B5}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.HasTable) {
            $tbl = $shp.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) Theme colours: Integral/Red Violet -> Office.
# ---------------------------------------------------------------------
function Get-RGBValue([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Index order for ThemeColorScheme.Colors(): dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink.
$officeColors = @(
    (Get-RGBValue 0x00 0x00 0x00), # dk1
    (Get-RGBValue 0xFF 0xFF 0xFF), # lt1
    (Get-RGBValue 0x44 0x54 0x6A), # dk2
    (Get-RGBValue 0xE7 0xE6 0xE6), # lt2
    (Get-RGBValue 0x5B 0x9B 0xD5), # accent1
    (Get-RGBValue 0xED 0x7D 0x31), # accent2
    (Get-RGBValue 0xA5 0xA5 0xA5), # accent3
    (Get-RGBValue 0xFF 0xC0 0x00), # accent4
    (Get-RGBValue 0x44 0x72 0xC4), # accent5
    (Get-RGBValue 0x70 0xAD 0x47), # accent6
    (Get-RGBValue 0x05 0x63 0xC1), # hlink
    (Get-RGBValue 0x95 0x4F 0x72)  # folHlink
)

$themeColorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($c = 1; $c -le $themeColorScheme.Count; $c++) {
    $themeColorScheme.Colors($c).RGB = $officeColors[$c - 1]
}
